# chore: update Sheets via scheduled runner
# Refresh market-board derived columns (H:N) for the affected leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2957.4443
$ws.Range("I76").Value = 2702.125
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 2702.125
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -2387.125
$ws.Range("N76").Value = -5630

$ws.Range("H79").Value = 2957.4443
$ws.Range("I79").Value = 2702.125
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 2702.125
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -1610.125
$ws.Range("N79").Value = -7184

$ws.Range("H113").Value = 620908.4
$ws.Range("I113").Value = 2224798.2
$ws.Range("J113").Value = 4027.6924
$ws.Range("K113").Value = 2224798.2
$ws.Range("L113").Value = 4027.6924
$ws.Range("M113").Value = -2221544.2
$ws.Range("N113").Value = -10535.6924

$ws.Range("H121").Value = 1271.24
$ws.Range("I121").Value = 140
$ws.Range("J121").Value = 1425.5
$ws.Range("K121").Value = 420
$ws.Range("L121").Value = 4276.5
$ws.Range("M121").Value = 1327
$ws.Range("N121").Value = -7770.5

$ws.Range("H137").Value = 2115.0476
$ws.Range("I137").Value = 1238.75
$ws.Range("J137").Value = 2654.3076
$ws.Range("K137").Value = 3716.25
$ws.Range("L137").Value = 7962.9228
$ws.Range("M137").Value = -1166.25
$ws.Range("N137").Value = -13062.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 7467.75
$ws.Range("I31").Value = 1706
$ws.Range("J31").Value = 47800
$ws.Range("K31").Value = 1706
$ws.Range("L31").Value = 47800
$ws.Range("M31").Value = -1412
$ws.Range("N31").Value = -48388

$ws.Range("H32").Value = 5230.89
$ws.Range("I32").Value = 4611.4634
$ws.Range("J32").Value = 17000
$ws.Range("K32").Value = 4611.4634
$ws.Range("L32").Value = 17000
$ws.Range("M32").Value = -4324.4634
$ws.Range("N32").Value = -17574

$ws.Range("H74").Value = 11906153
$ws.Range("I74").Value = 16130248
$ws.Range("J74").Value = 1883.5454
$ws.Range("K74").Value = 16130248
$ws.Range("L74").Value = 1883.5454
$ws.Range("M74").Value = -16129374
$ws.Range("N74").Value = -3631.5454

$ws.Range("H77").Value = 11906153
$ws.Range("I77").Value = 16130248
$ws.Range("J77").Value = 1883.5454
$ws.Range("K77").Value = 80651240
$ws.Range("L77").Value = 9417.726999999999
$ws.Range("M77").Value = -80646872
$ws.Range("N77").Value = -18153.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7814975
$ws.Range("I134").Value = 14707497
$ws.Range("J134").Value = 3451
$ws.Range("K134").Value = 44122491
$ws.Range("L134").Value = 10353
$ws.Range("M134").Value = -44119956
$ws.Range("N134").Value = -15423

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9437045
$ws.Range("I31").Value = 17242452
$ws.Range("J31").Value = 5513
$ws.Range("K31").Value = 17242452
$ws.Range("L31").Value = 5513
$ws.Range("M31").Value = -17242157
$ws.Range("N31").Value = -6103

$ws.Range("H34").Value = 9437045
$ws.Range("I34").Value = 17242452
$ws.Range("J34").Value = 5513
$ws.Range("K34").Value = 17242452
$ws.Range("L34").Value = 5513
$ws.Range("M34").Value = -17242250
$ws.Range("N34").Value = -5917

$ws.Range("H99").Value = 10204.286
$ws.Range("I99").Value = 14112.223
$ws.Range("J99").Value = 3170
$ws.Range("K99").Value = 14112.223
$ws.Range("L99").Value = 3170
$ws.Range("M99").Value = -12614.223
$ws.Range("N99").Value = -6166

$ws.Range("H126").Value = 10204.286
$ws.Range("I126").Value = 14112.223
$ws.Range("J126").Value = 3170
$ws.Range("K126").Value = 42336.669
$ws.Range("L126").Value = 9510
$ws.Range("M126").Value = -39866.669
$ws.Range("N126").Value = -14450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 10462.162
$ws.Range("I87").Value = 1225
$ws.Range("J87").Value = 11581.818
$ws.Range("K87").Value = 3675
$ws.Range("L87").Value = 34745.454
$ws.Range("M87").Value = -2427
$ws.Range("N87").Value = -37241.454

$ws.Range("H90").Value = 10462.162
$ws.Range("I90").Value = 1225
$ws.Range("J90").Value = 11581.818
$ws.Range("K90").Value = 11025
$ws.Range("L90").Value = 104236.362
$ws.Range("M90").Value = -4785
$ws.Range("N90").Value = -116716.362

$ws.Range("H113").Value = 467360.7
$ws.Range("I113").Value = 1221477.9
$ws.Range("J113").Value = 526.2381
$ws.Range("K113").Value = 3664433.7
$ws.Range("L113").Value = 1578.7143
$ws.Range("M113").Value = -3662263.7
$ws.Range("N113").Value = -5918.7143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1793.9445
$ws.Range("I97").Value = 1407.5
$ws.Range("J97").Value = 2566.8333
$ws.Range("K97").Value = 1407.5
$ws.Range("L97").Value = 2566.8333
$ws.Range("M97").Value = -911.5
$ws.Range("N97").Value = -3558.8333

$ws.Range("H102").Value = 3138.1177
$ws.Range("I102").Value = 2395.3635
$ws.Range("J102").Value = 4499.8335
$ws.Range("K102").Value = 2395.3635
$ws.Range("L102").Value = 4499.8335
$ws.Range("M102").Value = -773.3634999999999
$ws.Range("N102").Value = -7743.8335

$ws.Range("H122").Value = 2733.0667
$ws.Range("I122").Value = 1990
$ws.Range("J122").Value = 2847.3845
$ws.Range("K122").Value = 5970
$ws.Range("L122").Value = 8542.1535
$ws.Range("M122").Value = -3520
$ws.Range("N122").Value = -13442.1535

$ws.Range("H126").Value = 2219.2307
$ws.Range("I126").Value = 2031.25
$ws.Range("J126").Value = 2520
$ws.Range("K126").Value = 6093.75
$ws.Range("L126").Value = 7560
$ws.Range("M126").Value = -3623.75
$ws.Range("N126").Value = -12500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 433.42426
$ws.Range("I22").Value = 228.125
$ws.Range("J22").Value = 626.64703
$ws.Range("K22").Value = 228.125
$ws.Range("L22").Value = 626.64703
$ws.Range("M22").Value = 66.875
$ws.Range("N22").Value = -1216.64703

$ws.Range("H27").Value = 433.42426
$ws.Range("I27").Value = 228.125
$ws.Range("J27").Value = 626.64703
$ws.Range("K27").Value = 228.125
$ws.Range("L27").Value = 626.64703
$ws.Range("M27").Value = -121.125
$ws.Range("N27").Value = -840.64703

$ws.Range("H40").Value = 2681.2727
$ws.Range("I40").Value = 1972.8889
$ws.Range("J40").Value = 3531.3333
$ws.Range("K40").Value = 1972.8889
$ws.Range("L40").Value = 3531.3333
$ws.Range("M40").Value = -1836.8889
$ws.Range("N40").Value = -3803.3333

$ws.Range("H122").Value = 8046
$ws.Range("I122").Value = 9265.412
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 27796.236
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -25346.236
$ws.Range("N122").Value = -16600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2674592.8
$ws.Range("I113").Value = 6536287.5
$ws.Range("J113").Value = 1111.6923
$ws.Range("K113").Value = 19608862.5
$ws.Range("L113").Value = 3335.0769
$ws.Range("M113").Value = -19606692.5
$ws.Range("N113").Value = -7675.0769

$ws.Range("H122").Value = 1762.7059
$ws.Range("I122").Value = 1156.6
$ws.Range("J122").Value = 2628.5715
$ws.Range("K122").Value = 3469.8
$ws.Range("L122").Value = 7885.7145
$ws.Range("M122").Value = -1019.8
$ws.Range("N122").Value = -12785.7145

$ws.Range("H126").Value = 1726.5883
$ws.Range("I126").Value = 806.6
$ws.Range("J126").Value = 2109.9167
$ws.Range("K126").Value = 2419.8
$ws.Range("L126").Value = 6329.750100000001
$ws.Range("M126").Value = 50.19999999999982
$ws.Range("N126").Value = -11269.7501

$ws.Range("H136").Value = 4944.4863
$ws.Range("I136").Value = 1366.5
$ws.Range("J136").Value = 10192.2
$ws.Range("K136").Value = 4099.5
$ws.Range("L136").Value = 30576.6
$ws.Range("M136").Value = -1549.5
$ws.Range("N136").Value = -35676.60000000001
